$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move B4's value+style ("B4" text, green-bold style) over to E3.
$ws.Range("B4").Copy($ws.Range("E3"))

# B4 now becomes a plain cell holding the text that used to live in C1,
# with no special formatting.
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "C1"

# C1 is cleared out entirely, but keep the (now blank) cell record in
# place by touching its formatting once more after clearing.
$ws.Range("C1").ClearContents()
$ws.Range("C1").Font.Bold = $false
